$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'64.376.17"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "'3.493.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.52%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'586.03"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.09%  "
$ws.Range("D6").Value = "'135.01"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +1.92%  "
$ws.Range("D7").Value = "'3.493.91"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.53%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -1.09%  "
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "'7.20"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.27%  "
$ws.Range("E12").Value = "  -3.35%  "
$ws.Range("D13").Value = "'4.092.53"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.50%  "
$ws.Range("D14").Value = "'0.0000180"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +0.19%  "
$ws.Range("E15").Value = "  +1.15%  "
$ws.Range("D16").Value = "'3.499.02"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.38%  "
$ws.Range("D17").Value = "'64.381.97"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.45%  "
$ws.Range("D18").Value = "'25.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -9.98%  "
$ws.Range("D19").Value = "'10.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.62%  "
$ws.Range("E20").Value = "  -0.65%  "
$ws.Range("D21").Value = "'13.74"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.03%  "
$ws.Range("D22").Value = "'384.83"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.66%  "
$ws.Range("D23").Value = "'0.567"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.96%  "
$ws.Range("D24").Value = "'3.638.17"
$ws.Range("D24").Style = "Normal"
$ws.Range("D25").Value = "'74.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.48%  "
$ws.Range("D26").Value = "'1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.01%  "
$ws.Range("E27").Value = "  +3.23%  "
$ws.Range("B28").Value = "Fetch.AI"
$ws.Range("C28").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D28").Value = "'1.55"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.79%  "
$ws.Range("B29").Value = "RenderToken"
$ws.Range("C29").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D29").Value = "'7.47"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "
$ws.Range("E30").Value = "  +0.10%  "
$ws.Range("E31").Value = "  -1.44%  "
$ws.Range("D32").Value = "'8.22"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("D33").Value = "'3.516.51"
$ws.Range("D33").Style = "Normal"
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").Value = "'0.147"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.57%  "
$ws.Range("E36").Value = "  -2.22%  "
$ws.Range("D37").Value = "'5.29"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.19%  "
$ws.Range("E38").Value = "  -3.42%  "
$ws.Range("D39").Value = "'6.83"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.35%  "
$ws.Range("D40").Value = "'162.08"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.23%  "
$ws.Range("D41").Value = "'0.0781"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -3.19%  "
$ws.Range("D42").Value = "'0.804"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -1.24%  "
$ws.Range("D43").Value = "'25.76"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -3.18%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "'41.88"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.59%  "
$ws.Range("E46").Value = "  -0.03%  "
$ws.Range("D47").Value = "'4.42"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.48%  "
$ws.Range("D48").Value = "'1.64"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.93%  "
$ws.Range("D49").Value = "'2.470.98"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.98%  "
$ws.Range("D50").Value = "'6.77"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.88%  "
$ws.Range("D51").Value = "'0.907"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.67%  "
